# Insert a new worksheet "testAddItemToCart" as the first sheet in the
# workbook, containing the expected-message data used by an "add item to
# cart" test, and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Add the new sheet immediately before the current first sheet so it ends
# up in position 1, then name it.
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$newSheet.Name = "testAddItemToCart"

# Populate the two data cells.
$newSheet.Range("A1").Value = "expected_message"
$newSheet.Range("A2").Value = "Product successfully added to your shopping cart"

# Match the saved selection/active cell on the new sheet.
$newSheet.Range("A2").Select()

# Ensure the new sheet (now the first tab) is the active sheet.
$newSheet.Activate()
